$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 324
$ws1.Range("F5").Value = 345
$ws1.Range("F9").Value = 2268
$ws1.Range("F15").Value = 2184
$ws1.Range("F16").Value = 321

# Sheet "演出" (Shows)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 97

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 2095

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 2095
$ws4.Range("F9").Value = 324
$ws4.Range("F15").Value = 345
$ws4.Range("F19").Value = 97
$ws4.Range("F23").Value = 2268
$ws4.Range("F31").Value = 2184
$ws4.Range("F32").Value = 321
